$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A72").NumberFormat = "@"
$ws.Range("A72").Value = "2025/10/07"
$ws.Range("A72").ClearFormats()
$ws.Range("B72").Value = "火"
$ws.Range("C72").Value = 1
$ws.Range("D72").Value = 201
